# Actualización automática 2025-09-26 08:30:09
$wb = $excel.ActiveWorkbook

# --- Hoja "VENTAS POR GRUPO" ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("M4").Value = 347.95
$wsGrupo.Range("M38").Value = 4901.88
$wsGrupo.Range("D45").Value = 950.4
$wsGrupo.Range("M45").Value = 2252.16

# --- Hoja "VENTA MENSUAL" ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F4").Value = 794.99
$wsMensual.Range("F38").Value = 8278.48
$wsMensual.Range("F45").Value = 4109.99
$wsMensual.Range("F57").Value = 69497.67999999999

# --- Hoja "CUMPLIMIENTO MENSUAL" ---
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$wsCumplimiento.Range("D3").Value = 15402.15
$wsCumplimiento.Range("E3").Value = 4985.327421713499
$wsCumplimiento.Range("F3").Value = 0.7554711002939519

$wsCumplimiento.Range("D11").Value = 6200.76
$wsCumplimiento.Range("E11").Value = 13372.3002492497
$wsCumplimiento.Range("F11").Value = 0.3168007414802546

$wsCumplimiento.Range("D12").Value = 44144.29
$wsCumplimiento.Range("E12").Value = 2989.973157909801
$wsCumplimiento.Range("F12").Value = 0.9365647629221919

$wsCumplimiento.Range("D15").Value = 87975.45
$wsCumplimiento.Range("E15").Value = 10922.54992509275
$wsCumplimiento.Range("F15").Value = 0.8895574234730157
